$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.177.32'
$ws.Range("E2").Value = '  +2.11%  '

$ws.Range("D3").Value = '3.470.70'
$ws.Range("E3").Value = '  +1.74%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.09'
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.60'
$ws.Range("E6").Value = '  +3.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +10.30%  '

$ws.Range("D9").Value = '3.476.62'
$ws.Range("E9").Value = '  +1.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.25'
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +2.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +4.23%  '

$ws.Range("D13").Value = '4.057.97'
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("E14").Value = '  +0.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000193'
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.56'
$ws.Range("E16").Value = '  +3.80%  '

$ws.Range("D17").Value = '65.140.99'
$ws.Range("E17").Value = '  +2.08%  '

$ws.Range("D18").Value = '3.481.39'
$ws.Range("E18").Value = '  +2.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.51'
$ws.Range("E19").Value = '  +3.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.44'
$ws.Range("E20").Value = '  +2.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.16'
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.16'
$ws.Range("E22").Value = '  +1.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").Value = '  +4.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.15'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("E26").Value = '  +0.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.06'
$ws.Range("E27").Value = '  +6.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.52'
$ws.Range("E30").Value = '  +9.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.24'
$ws.Range("E31").Value = '  +1.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("E32").Value = '  +2.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.65'
$ws.Range("E33").Value = '  +1.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.31'
$ws.Range("E34").Value = '  +6.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +11.91%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.69'
$ws.Range("E36").Value = '  +1.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.94'
$ws.Range("E37").Value = '  +5.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0783'
$ws.Range("E38").Value = '  +3.26%  '

$ws.Range("D39").Value = '2.941.70'
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.83'
$ws.Range("E40").Value = '  +5.64%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.71'
$ws.Range("E41").Value = '  +8.52%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.86'
$ws.Range("E42").Value = '  +0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0324'
$ws.Range("E43").Value = '  +2.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.03'
$ws.Range("E44").Value = '  +2.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.783'
$ws.Range("E45").Value = '  +3.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.08'
$ws.Range("E46").Value = '  +12.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  +3.25%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '321.62'
$ws.Range("E48").Value = '  +9.93%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.112'
$ws.Range("E49").Value = '  +6.81%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.883'
$ws.Range("E50").Value = '  +5.64%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.19'
$ws.Range("E51").Value = '  -0.38%  '
